# The workbook tracks, per row, a Python-set-literal string in column F
# ("possible 84000 IDs"). This edit re-serializes several of those set
# literals with a different (but set-equivalent) element order — the same
# transformation applied uniformly everywhere that exact literal occurs.
# Use whole-cell text Replace over the used range so every matching row
# picks up the new ordering in one shot per distinct literal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.UsedRange

$pairs = @(
    @("{'eft:jnanasidhi', 'eft:jnanasiddhi'}", "{'eft:jnanasiddhi', 'eft:jnanasidhi'}"),
    @("{'eft:ch-nyi-tsultrim', 'eft:dharmatasila'}", "{'eft:dharmatasila', 'eft:ch-nyi-tsultrim'}"),
    @("{'eft:vidyakarasimha', 'eft:t-vidyakarasimha'}", "{'eft:t-vidyakarasimha', 'eft:vidyakarasimha'}"),
    @("{'eft:ban-de-dpal-brtsegs', 'eft:kawa-paltsek-under-the-name-paltsek-raksita-', 'eft:ska-ba-dpal-brtsegs', 'eft:dpal-brtsegs', 'eft:paltsek'}", "{'eft:dpal-brtsegs', 'eft:paltsek', 'eft:ska-ba-dpal-brtsegs', 'eft:ban-de-dpal-brtsegs', 'eft:kawa-paltsek-under-the-name-paltsek-raksita-'}"),
    @("{'eft:srilendrabodhi', 'eft:surendrabodhi', 'eft:silendrabodhi'}", "{'eft:surendrabodhi', 'eft:silendrabodhi', 'eft:srilendrabodhi'}"),
    @("{'eft:jinamitra', 'eft:jinamitra-k-', 'eft:dzi-na-mi-tra-k-'}", "{'eft:jinamitra', 'eft:dzi-na-mi-tra-k-', 'eft:jinamitra-k-'}"),
    @("{'eft:band-yesh-d-', 'eft:zhang-yesh-d-', 'eft:ye-shes-sde', 'eft:band-yesh-de', 'eft:yesh-d-ye-shes-sde-', 'eft:yesh-d-'}", "{'eft:yesh-d-', 'eft:yesh-d-ye-shes-sde-', 'eft:ye-shes-sde', 'eft:band-yesh-de', 'eft:zhang-yesh-d-', 'eft:band-yesh-d-'}"),
    @("{'eft:munivarma', 'eft:munivarman'}", "{'eft:munivarman', 'eft:munivarma'}"),
    @("{'eft:ban-de-dpal-gyi-lhun-po', 'eft:palgyi-lh-npo', 'eft:dpal-gyi-lhun-po'}", "{'eft:ban-de-dpal-gyi-lhun-po', 'eft:dpal-gyi-lhun-po', 'eft:palgyi-lh-npo'}"),
    @("{'eft:dipamkarasrijnana', 'eft:dipamkara-srijnana'}", "{'eft:dipamkara-srijnana', 'eft:dipamkarasrijnana'}")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $rng.Replace($old, $new, 1, 1, $false, $false, $false, $false) | Out-Null
}
